$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the English column headers with their Persian translations
# (feature: support Persian titles on Excel import).
$ws.Range("A1").Value = "ستون اول"
$ws.Range("B1").Value = "ستون دوم"
$ws.Range("C1").Value = "ستون سوم"

# Leave the cursor/selection on C2, matching where the author left it
$ws.Range("C2").Select()
